# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Thu Jan 18 20:26:42 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '40.957.98'
$ws.Range('E2').Value = '  -4.10%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.441.76'
$ws.Range('E3').Value = '  -3.59%  '

# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.05%  '

# Row 5: BNB
$ws.Range('D5').Value = '''309.47'
$ws.Range('E5').Value = '  +0.30%  '

# Row 6: Solana
$ws.Range('D6').Value = '''93.12'
$ws.Range('E6').Value = '  -7.61%  '

# Row 7: XRP
$ws.Range('E7').Value = '  -3.66%  '

# Row 8: USDC
$ws.Range('E8').Value = '  +0.13%  '

# Row 9: Cardano
$ws.Range('D9').Value = '''0.501'
$ws.Range('E9').Value = '  -5.15%  '

# Row 10: Avalanche
$ws.Range('D10').Value = '''33.21'
$ws.Range('E10').Value = '  -8.27%  '

# Row 11: Dogecoin
$ws.Range('E11').Value = '  -3.12%  '

# Row 12: TRON
$ws.Range('E12').Value = '  -0.69%  '

# Row 13: Polkadot
$ws.Range('E13').Value = '  -5.35%  '

# Row 14: WrappedliquidstakedEther2.0
$ws.Range('D14').Value = '2.822.63'
$ws.Range('E14').Value = '  -3.68%  '

# Row 15: WrappedEther
$ws.Range('D15').Value = '2.452.68'
$ws.Range('E15').Value = '  -4.26%  '

# Row 16: Chainlink
$ws.Range('D16').Value = '''14.37'
$ws.Range('E16').Value = '  -9.36%  '

# Row 17: Polygon
$ws.Range('E17').Value = '  -3.43%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '41.017.90'
$ws.Range('E18').Value = '  -3.89%  '

# Row 19: Uniswap
$ws.Range('E19').Value = '  -6.80%  '

# Row 20: ShibaInu
$ws.Range('D20').Value = '0.0₃0912'
$ws.Range('E20').Value = '  -4.37%  '

# Row 21: InternetComputer(DFINITY)
$ws.Range('D21').Value = '''11.43'
$ws.Range('E21').Value = '  -6.55%  '

# Row 22: Litecoin
$ws.Range('D22').Value = '''67.14'
$ws.Range('E22').Value = '  -3.26%  '

# Row 23: BitcoinCash
$ws.Range('D23').Value = '''236.28'
$ws.Range('E23').Value = '  -3.14%  '

# Row 24: PancakeSwap
$ws.Range('E24').Value = '  -4.38%  '

# Row 25: ImmutableX
$ws.Range('D25').Value = '''1.92'
$ws.Range('E25').Value = '  -5.90%  '

# Row 26: Dai
$ws.Range('E26').Value = '  +0.17%  '

# Row 27: EthereumClassic
$ws.Range('D27').Value = '''24.33'
$ws.Range('E27').Value = '  -6.54%  '

# Row 28: Toncoin
$ws.Range('D28').Value = '''2.22'
$ws.Range('E28').Value = '  -4.32%  '

# Row 29: Cosmos
$ws.Range('E29').Value = '  -5.08%  '

# Row 30: InjectiveProtocol
$ws.Range('D30').Value = '''35.83'
$ws.Range('E30').Value = '  -8.78%  '

# Row 31: Monero
$ws.Range('D31').Value = '''152.07'
$ws.Range('E31').Value = '  -2.24%  '

# Row 32: Filecoin
$ws.Range('D32').Value = '''5.57'
$ws.Range('E32').Value = '  -3.75%  '

# Row 33: WEMIXToken
$ws.Range('D33').Value = '''2.60'
$ws.Range('E33').Value = '  -0.83%  '

# Row 34: Hedera
$ws.Range('D34').Value = '''0.0749'
$ws.Range('E34').Value = '  -5.47%  '

# Row 35: ApeXProtocol
$ws.Range('E35').Value = '  -8.51%  '

# Row 36: LidoDAOToken
$ws.Range('D36').Value = '''2.99'
$ws.Range('E36').Value = '  -5.45%  '

# Row 37: Celestia
$ws.Range('E37').Value = '  -6.94%  '

# Row 38: ARBITRUM
$ws.Range('E38').Value = '  -7.73%  '

# Row 39: Kaspa
$ws.Range('D39').Value = '''0.103'
$ws.Range('E39').Value = '  -8.45%  '

# Row 40: Stellar
$ws.Range('E40').Value = '  -4.36%  '

# Row 41: RenderToken
$ws.Range('D41').Value = '''4.10'
$ws.Range('E41').Value = '  -4.98%  '

# Row 42: EnergySwap
$ws.Range('D42').Value = '''21.05'
$ws.Range('E42').Value = '  -4.62%  '

# Row 43: FirstDigitalUSD
$ws.Range('E43').Value = '  +0.06%  '

# Row 44: Maker
$ws.Range('D44').Value = '1.969.05'

# Row 45: VeChain
$ws.Range('E45').Value = '  -5.18%  '

# Row 46: NEARProtocol
$ws.Range('E46').Value = '  -8.20%  '

# Row 47: FraxShare
$ws.Range('D47').Value = '''8.68'
$ws.Range('E47').Value = '  -2.40%  '

# Row 48: BitcoinSV
$ws.Range('D48').Value = '''76.74'
$ws.Range('E48').Value = '  -4.95%  '

# Row 49: Aave
$ws.Range('D49').Value = '''96.72'
$ws.Range('E49').Value = '  -4.22%  '

# Row 50: ordi -> Algorand
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '''0.179'
$ws.Range('E50').Value = '  -6.88%  '

# Row 51: Algorand -> ordi
$ws.Range('B51').Value = 'ordi'
$ws.Range('C51').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D51').Value = '''68.32'
$ws.Range('E51').Value = '  -5.89%  '
